$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking crypto price/volume refresh (GitHub Actions bot).
# D-column "Price" cells are plain text in the source sheet; Excel's
# COM layer auto-coerces bare numeric-looking text to a Number on
# assignment, so we briefly force a Text format, write the literal
# string, then restore the Normal style (clears the number-format
# override again so the cell carries no stray style index).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.864.30'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.18'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.65%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  +1.94%  '
$ws.Range("E7").Value = '  -0.44%  '
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0622'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.47%  '
$ws.Range("E11").Value = '  +0.44%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.43'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.63%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.659.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.81%  '
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.528'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.44%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.25'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.70%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.864.08'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.88%  '
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").Value = '  +1.87%  '
$ws.Range("E22").Value = '  +7.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.18%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.39'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.75%  '
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0506'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.13%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("E32").Value = '  -0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("E34").Value = '  +2.55%  '
$ws.Range("E35").Value = '  -0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.244.62'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.96%  '
$ws.Range("E37").Value = '  +1.07%  '
$ws.Range("E38").Value = '  +2.93%  '
$ws.Range("E39").Value = '  +3.94%  '
$ws.Range("E40").Value = '  -0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.805'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.95%  '
$ws.Range("E42").Value = '  +2.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.783.61'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("E44").Value = '  -2.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.74'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '91.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("E48").Value = '  +16.90%  '
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("E50").Value = '  +2.29%  '
$ws.Range("E51").Value = '  +2.13%  '
